# "Add files via upload" — appends a new "Ryan" / "63.35.235.192" record
# to the VM Allocations list (columns: A = IP Address, B = Name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18. Write the Name (column B) before the IP Address (column A)
# so new shared-string entries land in the same order as the source edit
# (Ryan, then the IP).
$ws.Range("B18").Value = "Ryan"
$ws.Range("A18").Value = "63.35.235.192"

# The IP-address cell (A18) carries a distinct font: Ubuntu 12pt, RGB
# #212529 (21,25,29 -> 0x292521 in BGR/VBA RGB() encoding).
$ws.Range("A18").Font.Name = "Ubuntu"
$ws.Range("A18").Font.Size = 12
$ws.Range("A18").Font.Color = 2696481
$ws.Range("A18").RowHeight = 16.5

# Match the saved selection left behind in the workbook.
$ws.Range("C22").Select()
